# Update "想去人数" (F) and "最低票价" (G) numbers on the "展览" and
# "全部类型" worksheets to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - full set of updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 3390
$ws1.Range("G3").Value  = 75
$ws1.Range("G4").Value  = 65
$ws1.Range("F7").Value  = 1739
$ws1.Range("F14").Value = 39
$ws1.Range("F20").Value = 12
$ws1.Range("F24").Value = 46
$ws1.Range("F26").Value = 404
$ws1.Range("F27").Value = 265
$ws1.Range("F28").Value = 117
$ws1.Range("F29").Value = 44
$ws1.Range("F32").Value = 447
$ws1.Range("F33").Value = 2314
$ws1.Range("F37").Value = 570
$ws1.Range("F38").Value = 565
$ws1.Range("F40").Value = 239
$ws1.Range("F41").Value = 357
$ws1.Range("F43").Value = 544

# Sheet "全部类型" (all types) - same updates; F40 already matches the
# new value on this sheet, so it is left untouched.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 3390
$ws4.Range("G3").Value  = 75
$ws4.Range("G4").Value  = 65
$ws4.Range("F7").Value  = 1739
$ws4.Range("F14").Value = 39
$ws4.Range("F20").Value = 12
$ws4.Range("F24").Value = 46
$ws4.Range("F26").Value = 404
$ws4.Range("F27").Value = 265
$ws4.Range("F28").Value = 117
$ws4.Range("F29").Value = 44
$ws4.Range("F32").Value = 447
$ws4.Range("F33").Value = 2314
$ws4.Range("F37").Value = 570
$ws4.Range("F38").Value = 565
$ws4.Range("F41").Value = 357
$ws4.Range("F43").Value = 544
